# Add a new "2022" column (column N) to the table, mirroring the existing
# "2021" column (column M) both in values and in cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values for the new column N --------------------------------------
$ws.Range("N2").Value  = $null
$ws.Range("N3").Value  = 2022
$ws.Range("N4").Value  = 1434
$ws.Range("N5").Value  = 12822
$ws.Range("N6").Value  = 3099
$ws.Range("N7").Value  = 9722
$ws.Range("N8").Value  = 14424
$ws.Range("N9").Value  = 5279
$ws.Range("N10").Value = 9145

# --- Formatting: copy the style of each M-column cell onto the matching
#     N-column cell, row by row, so borders/fonts/number formats line up
#     exactly like the rest of the table. -------------------------------
$rows = 2..10
foreach ($r in $rows) {
    $ws.Range("M$r").Copy() | Out-Null
    $ws.Range("N$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

# --- Update the selection to match what was left after editing --------
$ws.Range("N2").Select() | Out-Null

$wb.Save()
